$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "col3_hum_chung_1974" data row (originally row 5). Deleting the
# entire row shifts every subsequent row up by one, which reproduces the
# renumbering seen throughout the diff (old row 6 -> new row 5, etc.) while
# keeping the rest of the sheet's row-to-label associations intact.
$ws.Range("A5:G5").EntireRow.Delete()

# After the delete, the sheet has data rows 2-47 (A2:G47). Append six new
# rows of data (Collagen / skin-collagen examples) at rows 48-53.

# Row 48: Collagen_1_hum_uitto_1978
$ws.Range("A48").Value = "Collagen_1_hum_uitto_1978"
$ws.Range("B48").Value = 0.1843275692130446
$ws.Range("C48").Value = 5.425124436183536
$ws.Range("E48").Value = 0.1169947657113922
$ws.Range("F48").Value = 0.1356288666985065
$ws.Range("G48").Value = 1.034308447330943

# Row 49: Collagen_3_hum_chung_1974 (reuses the values originally stored
# against col3_hum_chung_1974 before its row was removed above)
$ws.Range("A49").Value = "Collagen_3_hum_chung_1974"
$ws.Range("B49").Value = 0.1826407971543893
$ws.Range("C49").Value = 5.475227964290386
$ws.Range("E49").Value = 0.1555862244343698
$ws.Range("F49").Value = 0.1803669007380068
$ws.Range("G49").Value = 1.375481588789803

# Row 50: Collagen_4_hum_glanville_1979
$ws.Range("A50").Value = "Collagen_4_hum_glanville_1979"
$ws.Range("B50").Value = 0.1643559632507289
$ws.Range("C50").Value = 6.084354837034276
$ws.Range("E50").Value = 0.1695406238034565
$ws.Range("F50").Value = 0.1965438584025607
$ws.Range("G50").Value = 1.498847391158093

# Row 51: Skin-Collagen_hum_bornstein_1964
$ws.Range("A51").Value = "Skin-Collagen_hum_bornstein_1964"
$ws.Range("B51").Value = 0.1833505892952129
$ws.Range("C51").Value = 5.454032102345193
$ws.Range("E51").Value = 0.116479586664021
$ws.Range("F51").Value = 0.135031633566606
$ws.Range("G51").Value = 1.029753935534248

# Row 52: Skin-Collagen_hum_acid_miyahara_1978
$ws.Range("A52").Value = "Skin-Collagen_hum_acid_miyahara_1978"
$ws.Range("B52").Value = 0.183701390849173
$ws.Range("C52").Value = 5.443616922971718
$ws.Range("E52").Value = 0.114902561622547
$ws.Range("F52").Value = 0.1332034311010537
$ws.Range("G52").Value = 1.015812027004146

# Row 53: Skin-Collagen_hum_age0_miyahara_1978
$ws.Range("A53").Value = "Skin-Collagen_hum_age0_miyahara_1978"
$ws.Range("B53").Value = 0.1798439066934646
$ws.Range("C53").Value = 5.560377431660513
$ws.Range("E53").Value = 0.1148732530857138
$ws.Range("F53").Value = 0.1331694545072196
$ws.Range("G53").Value = 1.01555292082071

# Match the bold/bordered/centered label style used by the other rows in
# column A (copy formatting only, so the text values set above are kept).
$ws.Range("A2").Copy()
$ws.Range("A48:A53").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A1").Select()
